{"js": "// The document stores short \"<tag>value</tag>\" markers as separate runs,\n// e.g. \"<id>\", \"p137r_1\", \"</id>\" \u2014 each with its own (slightly different)\n// run-level formatting. This edit merges those three runs back into a\n// single run per paragraph (\"<id>p137r_1</id>\", etc.) using the\n// formatting of the surrounding \"<id>\"/\"</id>\" runs (Courier New).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].load(\"text\");\n}\nawait context.sync();\n\n// Match paragraphs whose whole text is exactly \"<id>NAME</id>\" where NAME\n// looks like this document's \"p137r_N\" ids (not e.g. \"fig_p137r_1\", which\n// stays untouched).\nconst idRegex = /^<id>(p137r_\\d+)<\\/id>$/;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const match = idRegex.exec(para.text);\n  if (!match) continue;\n\n  // Re-inserting the same full text as a \"Replace\" collapses the\n  // paragraph's runs into a single run, inheriting the formatting of the\n  // first run that previously held \"<id>\" (Courier New / 7f6000 / 18pt).\n  para.insertText(match[0], \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# The document stores short \"<tag>value</tag>\" markers as separate runs,\n# e.g. \"<id>\", \"p137r_1\", \"</id>\" -- each with its own (slightly different)\n# run-level formatting. This edit merges those three runs back into a\n# single run per paragraph (\"<id>p137r_1</id>\", etc.), keeping the\n# Courier New formatting that the \"<id>\"/\"</id>\" runs already used.\n$d = $word.ActiveDocument\n$count = $d.Paragraphs.Count\n\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $r = $p.Range\n    # Exclude the trailing paragraph mark from the text comparison.\n    $r.MoveEnd(1, -1) | Out-Null\n    $t = $r.Text\n\n    # Only paragraphs whose whole text is \"<id>p137r_N</id>\" (not e.g.\n    # \"<id>fig_p137r_1</id>\", which stays untouched).\n    if ($t -match '^(<id>)(p137r_\\d+)</id>$') {\n        $idOpen = $matches[1]\n        $idName = $matches[2]\n\n        # Keep the opening \"<id>\" run (and its Courier New formatting)\n        # untouched; delete everything after it in the paragraph...\n        $idEnd = $r.Start + $idOpen.Length\n        $rest = $d.Range($idEnd, $r.End)\n        $rest.Delete()\n\n        # ...then retype the rest right after it. Word inherits the\n        # formatting of the preceding text, so the new run picks up the\n        # same Courier New / color / size as the \"<id>\" run, and merges\n        # with it into a single run.\n        $insPoint = $d.Range($idEnd, $idEnd)\n        $insPoint.InsertAfter(\"$idName</id>\")\n    }\n}\n"}
